$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 7935
$ws.Range("F5").Value = 2151
$ws.Range("F6").Value = 8681
$ws.Range("F7").Value = 14
$ws.Range("F10").Value = 5809
$ws.Range("F12").Value = 2808
$ws.Range("F13").Value = 1219
$ws.Range("F15").Value = 102
$ws.Range("F18").Value = 117
$ws.Range("F19").Value = 4011
$ws.Range("F26").Value = 5802
$ws.Range("F28").Value = 79
$ws.Range("F32").Value = 419
$ws.Range("F33").Value = 4293
$ws.Range("F34").Value = 1555
$ws.Range("F37").Value = 5718
$ws.Range("F38").Value = 83
$ws.Range("F41").Value = 31
$ws.Range("F42").Value = 3712
$ws.Range("F43").Value = 22
$ws.Range("F45").Value = 2363

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 164
$ws.Range("F4").Value = 20

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1379

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1379
$ws.Range("F4").Value = 7935
$ws.Range("F5").Value = 2151
$ws.Range("F6").Value = 8681
$ws.Range("F7").Value = 14
$ws.Range("F9").Value = 5809
$ws.Range("F11").Value = 2808
$ws.Range("F12").Value = 1219
$ws.Range("F14").Value = 102
$ws.Range("F16").Value = 164
$ws.Range("F18").Value = 20
$ws.Range("F19").Value = 117
$ws.Range("F20").Value = 4011
$ws.Range("F27").Value = 5802
$ws.Range("F29").Value = 79
$ws.Range("F32").Value = 420
$ws.Range("F35").Value = 1555
$ws.Range("F39").Value = 5718
$ws.Range("F40").Value = 83
$ws.Range("F43").Value = 3712
$ws.Range("F44").Value = 2363
